# Adds the "femaleClients" (screenClient / client6Week / client6Month) and
# "maleClients" (screenPartner / partner6Month) test-form entries to the
# "framework" test harness workbook: one 3-row block per form on the
# "survey" sheet (external-link launcher pattern already used by every
# other test form there) plus the matching lookup rows on the "choices"
# sheet.

$wb = $excel.ActiveWorkbook
$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# 1) "survey" sheet: append 5 new 3-row blocks (rows 74-88) after the
#    existing last block (rows 71-73, "eonasdan").
#    Each block is: [name row] / [hyperlink-style "open form" row] /
#    [blank "exit section" row] - exactly the pattern already used for
#    every other test form above it.
# ---------------------------------------------------------------------

# The last existing row (73, "exit section" for eonasdan) gets a tiny
# row-height tweak as a side effect of the new rows being appended below.
$survey.Rows.Item(73).RowHeight = 17

$forms = @(
    @{ Name = "screenClient";  Hash = "femaleClients/forms/screenClient";  Label = "Screen Female Client"; HRow = 62 },
    @{ Name = "client6Week";   Hash = "femaleClients/forms/client6Week";   Label = "Client 6 Week";         HRow = 62 },
    @{ Name = "client6Month";  Hash = "femaleClients/forms/client6Month";  Label = "Client 6 Month";        HRow = 58 },
    @{ Name = "screenPartner"; Hash = "maleClients/forms/screenPartner";   Label = "Screen Partner";        HRow = 58 },
    @{ Name = "partner6Month"; Hash = "maleClients/forms/partner6Month";   Label = "Partner 6 Month";       HRow = 58 }
)

$row = 74
foreach ($f in $forms) {

    # Row 1 of the block: form id in column A, plain/default style.
    $survey.Rows.Item($row).RowHeight = 12.75
    $survey.Cells.Item($row, 1).Value = $f.Name

    # Row 2 of the block: column A keeps the "blank but formatted" style
    # used throughout (copy format from the template block above), and
    # column B holds the quoted getHashString() formula-looking text
    # (apostrophe doubled so it is kept as literal text, matching the
    # quote-prefixed style already used for every other entry), plus the
    # "external_link" / "Open form" marker cells.
    $hashRow = $row + 1
    $survey.Rows.Item($hashRow).RowHeight = $f.HRow
    $survey.Cells.Item(72, 1).Copy() | Out-Null
    $survey.Cells.Item($hashRow, 1).PasteSpecial(-4122) | Out-Null
    $survey.Cells.Item(72, 2).Copy() | Out-Null
    $survey.Cells.Item($hashRow, 2).PasteSpecial(-4122) | Out-Null
    $survey.Cells.Item($hashRow, 2).Value = "''?' + opendatakit.getHashString('../config/tables/" + $f.Hash + "/',null)"
    $survey.Cells.Item($hashRow, 5).Value = "external_link"
    $survey.Cells.Item($hashRow, 7).Value = "Open form"

    # Row 3 of the block: "exit section" marker, same formatting as every
    # other block's closing row.
    $exitRow = $row + 2
    $survey.Rows.Item($exitRow).RowHeight = 12.75
    $survey.Cells.Item(73, 1).Copy() | Out-Null
    $survey.Cells.Item($exitRow, 1).PasteSpecial(-4122) | Out-Null
    $survey.Cells.Item(73, 2).Copy() | Out-Null
    $survey.Cells.Item($exitRow, 2).PasteSpecial(-4122) | Out-Null
    $survey.Cells.Item($exitRow, 3).Value = "exit section"

    $row += 3
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) "choices" sheet: append the matching lookup rows (24-28) under the
#    "test_forms" choice list, same 3-column layout as every row above.
# ---------------------------------------------------------------------

$choiceRows = @(
    @{ Value = "screenClient";  Label = "Screen Female Client" },
    @{ Value = "client6Week";   Label = "Client 6 Week" },
    @{ Value = "client6Month";  Label = "Client 6 Month" },
    @{ Value = "screenPartner"; Label = "Screen Partner" },
    @{ Value = "partner6Month"; Label = "Partner 6 Month" }
)

$crow = 24
foreach ($c in $choiceRows) {
    $choices.Cells.Item($crow, 1).Value = "test_forms"
    $choices.Cells.Item($crow, 2).Value = $c.Value
    $choices.Cells.Item($crow, 3).Value = $c.Label
    $crow += 1
}

# ---------------------------------------------------------------------
# 3) View-state: leave "choices" as the active sheet (unchanged), but
#    update the remembered selection on both touched sheets to reflect
#    where editing ended up.
# ---------------------------------------------------------------------

$survey.Activate()
$survey.Range("B80").Select()

$choices.Activate()
$choices.Range("H35").Select()

Write-Output "done"
